$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.143.43"
$ws.Range("E2").Value = "  -4.68%  "
$ws.Range("D3").Value = "1.656.26"
$ws.Range("E3").Value = "  -3.30%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").Value = "218.09"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").Value = "0.5148"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "0.06420"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").Value = "0.2560"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "19.85"
$ws.Range("E10").Value = "  -5.31%  "
$ws.Range("D11").Value = "0.07689"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.659.07"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "1.890.62"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "4.295"
$ws.Range("E14").Value = "  -6.06%  "
$ws.Range("D15").Value = "0.5528"
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").Value = "0.0₅8010"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "64.28"
$ws.Range("E17").Value = "  -5.52%  "
$ws.Range("D18").Value = "26.181.53"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "209.89"
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("D21").Value = "4.392"
$ws.Range("E21").Value = "  -6.14%  "
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").Value = "5.858"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "1.011"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "144.36"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "1.758"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "0.1158"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").Value = "6.945"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("D29").Value = "15.72"
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "0.05246"
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").Value = "1.259"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "3.359"
$ws.Range("E32").Value = "  -4.50%  "
$ws.Range("D33").Value = "3.220"
$ws.Range("E33").Value = "  -6.35%  "
$ws.Range("D34").Value = "1.562"
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.381"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "2.744"
$ws.Range("E36").Value = "  -4.61%  "
$ws.Range("D37").Value = "0.9226"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").Value = "0.5761"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "1.153.21"
$ws.Range("E39").Value = "  +10.11%  "
$ws.Range("D40").Value = "0.01590"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").Value = "1.011"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").Value = "0.8377"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "5.656"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D44").Value = "99.79"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "1.798.33"
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("E46").Value = "  -4.75%  "
$ws.Range("D47").Value = "0.4509"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "55.86"
$ws.Range("E48").Value = "  -3.87%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "7.876"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("D51").Value = "0.05087"
$ws.Range("E51").Value = "  -2.92%  "
